# Commit: "updated password for wordpress"
# Adds two new Wordpress Blog accounts (Angela, Laurie) as new rows below
# the existing data, each with a hyperlinked email address in column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 16: Wordpress Blog - Angela ---
$ws.Range("A16").Value = "Wordpress Blog - Angela"
$ws.Range("B16").Value = "angelal4"
$ws.Range("D16").Value = "angelgirl2272@gmail.com"
$ws.Hyperlinks.Add($ws.Range("D16"), "mailto:angelgirl2272@gmail.com")
$ws.Range("D16").Style = "Hyperlink"

# --- Row 17: Wordpress Blog - Laurie ---
$ws.Range("A17").Value = "Wordpress Blog - Laurie"
$ws.Range("B17").Value = "laurie415"
$ws.Range("D17").Value = "ldxtran@gmail.com"
$ws.Hyperlinks.Add($ws.Range("D17"), "mailto:ldxtran@gmail.com")
$ws.Range("D17").Style = "Hyperlink"

# Restore the selected cell as in the authored workbook
$ws.Range("D12").Select()
